$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($ws, $r, $data) {
    foreach ($col in $data.Keys) {
        $addr = "$col$r"
        $v = $data[$col]
        if ($v -eq $null) {
            $ws.Range($addr).Value = ""
        } else {
            $ws.Range($addr).Value = $v
        }
    }
}

$row2 = @{ "A" = 80469805; "B" = 90645; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 4361; "F" = "'Orange taggsvamp"; "G" = "'Hydnellum aurantiacum"; "H" = "'(Batsch:Fr.) P.Karst."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714440.0440998246; "R" = 6623501.154866487; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = "'Allmänt förekommande."; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row3 = @{ "A" = 80469685; "B" = 5113; "C" = "'Ovaliderad"; "D" = "'LC"; "E" = 100526; "F" = "'Bronshjon"; "G" = "'Callidium coriaceum"; "H" = "'Paykull, 1800"; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = "'äldre gnagspår"; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714517.2209692324; "R" = 6623428.185489144; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row4 = @{ "A" = 80469740; "B" = 90676; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 5966; "F" = "'Motaggsvamp"; "G" = "'Sarcodon squamosus"; "H" = "'(Schaeff.) Quél."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714373.2229025841; "R" = 6623438.118456176; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row5 = @{ "A" = 80469729; "B" = 56411; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 100049; "F" = "'Spillkråka"; "G" = "'Dryocopus martius"; "H" = "'(Linnaeus, 1758)"; "I" = "'1"; "J" = $null; "K" = $null; "L" = $null; "M" = "'födosökande"; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714421.0537466849; "R" = 6623434.797179379; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row6 = @{ "A" = 80469761; "B" = 90645; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 4361; "F" = "'Orange taggsvamp"; "G" = "'Hydnellum aurantiacum"; "H" = "'(Batsch:Fr.) P.Karst."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714390.1611933979; "R" = 6623460.847534781; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row7 = @{ "A" = 80469731; "B" = 103265; "C" = "'Ovaliderad"; "D" = "'LC"; "E" = 221144; "F" = "'Grönpyrola"; "G" = "'Pyrola chlorantha"; "H" = "'Sw."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714421.0537466849; "R" = 6623434.797179379; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row8 = @{ "A" = 80469819; "B" = 90661; "C" = "'Ovaliderad"; "D" = "'VU"; "E" = 2058; "F" = "'Koppartaggsvamp"; "G" = "'Hydnellum lundellii"; "H" = "'(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg"; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714455.1568480013; "R" = 6623493.929075443; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row9 = @{ "A" = 80469789; "B" = 90645; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 4361; "F" = "'Orange taggsvamp"; "G" = "'Hydnellum aurantiacum"; "H" = "'(Batsch:Fr.) P.Karst."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714409.7699023501; "R" = 6623481.20045886; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = "'Allmänt förekommande."; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row10 = @{ "A" = 80469792; "B" = 5135; "C" = "'Ovaliderad"; "D" = "'LC"; "E" = 105930; "F" = "'Vågbandad barkbock"; "G" = "'Semanotus undatus"; "H" = "'(Linnaeus, 1758)"; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = "'äldre gnagspår"; "N" = $null; "P" = "'Getängsjön, V om, Upl"; "Q" = 714409.7699023501; "R" = 6623481.20045886; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row11 = @{ "A" = 80469829; "B" = 90642; "C" = "'Ovaliderad"; "D" = "'VU"; "E" = 150; "F" = "'Grangråticka"; "G" = "'Boletopsis leucomelaena"; "H" = "'(Pers.) Fayod"; "I" = "'10"; "J" = "'fruktkroppar"; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714526.0508117813; "R" = 6623538.987286027; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row12 = @{ "A" = 80469832; "B" = 98520; "C" = "'Ovaliderad"; "D" = "'LC"; "E" = 222498; "F" = "'Blåsippa"; "G" = "'Hepatica nobilis"; "H" = "'Schreb."; "I" = $null; "J" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714526.0508117813; "R" = 6623538.987286027; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }
$row13 = @{ "A" = 80469695; "B" = 90676; "C" = "'Ovaliderad"; "D" = "'NT"; "E" = 5966; "F" = "'Motaggsvamp"; "G" = "'Sarcodon squamosus"; "H" = "'(Schaeff.) Quél."; "I" = "'8"; "J" = "'fruktkroppar"; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "P" = "'Getängsjön, N om, Upl"; "Q" = 714514.1731204849; "R" = 6623437.117110151; "S" = 5; "T" = "'Stockholm"; "U" = "'Norrtälje"; "V" = "'Uppland"; "W" = "'Frötuna"; "Y" = "'2019-10-09"; "Z" = "'00:00"; "AA" = "'2019-10-09"; "AB" = "'00:00"; "AC" = $null; "AD" = $false; "AE" = $false; "AF" = $null; "AG" = $false; "AT" = $null; "AW" = "'Bo Törnquist"; "AX" = "'Bo Törnquist, Kjell  Andersson"; "AY" = $null }

Set-RowData $ws 2 $row2
Set-RowData $ws 3 $row3
Set-RowData $ws 4 $row4
Set-RowData $ws 5 $row5
Set-RowData $ws 6 $row6
Set-RowData $ws 7 $row7
Set-RowData $ws 8 $row8
Set-RowData $ws 9 $row9
Set-RowData $ws 10 $row10
Set-RowData $ws 11 $row11
Set-RowData $ws 12 $row12
Set-RowData $ws 13 $row13
